$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6913
$ws1.Range("F3").Value = 17
$ws1.Range("F5").Value = 78
$ws1.Range("F8").Value = 121
$ws1.Range("F9").Value = 111
$ws1.Range("F11").Value = 10
$ws1.Range("F13").Value = 191
$ws1.Range("F14").Value = 435
$ws1.Range("F16").Value = 1807
$ws1.Range("F17").Value = 38
$ws1.Range("F18").Value = 3525
$ws1.Range("F19").Value = 24
$ws1.Range("F20").Value = 240
$ws1.Range("F22").Value = 2149
$ws1.Range("F23").Value = 211
$ws1.Range("F27").Value = 2
$ws1.Range("F29").Value = 144

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6913
$ws4.Range("F3").Value = 17
$ws4.Range("F5").Value = 78
$ws4.Range("F9").Value = 121
$ws4.Range("F10").Value = 111
$ws4.Range("F12").Value = 10
$ws4.Range("F14").Value = 191
$ws4.Range("F15").Value = 435
$ws4.Range("F17").Value = 1807
$ws4.Range("F18").Value = 38
$ws4.Range("F19").Value = 3525
$ws4.Range("F20").Value = 24
$ws4.Range("F21").Value = 240
$ws4.Range("F23").Value = 2149
$ws4.Range("F24").Value = 211
$ws4.Range("F28").Value = 2
$ws4.Range("F30").Value = 144
